$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set H2 value to the Chinese character meaning "Yes" (是)
$ws.Range("H2").Value = "是"

# Move the active selection to H12 (as shown by the new selection in the sheet view)
$ws.Range("H12").Select()
